# Actualización automática del mapa (2025-08-26 15:13:04)
# Refreshes rows 44 and 85-89 of the "AYKO" sheet with the latest feed data.
# Numeric-looking / date-looking strings are written with a leading
# apostrophe so Excel stores them as literal text (matching the source
# data, which keeps "Caso"/"OT"/"Comuna"/date columns as text) instead of
# silently coercing them into numbers or date serials.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 44 -----------------------------------------------------------
$ws.Range("A44").Value = "'7054"
$ws.Range("C44").Value = "HUMAHUACA 3918"

# --- Row 85 -------------------------------------------------------------
$ws.Range("A85").Value = "'6960"
$ws.Range("C85").Value = "VALLESE, FELIPE 1940"
$ws.Range("D85").Value = "'7"
$ws.Range("E85").Value = "'808972988"
$ws.Range("M85").Value = -58.460818
$ws.Range("N85").Value = -34.618934
$ws.Range("O85").Value = "Boedo"

# --- Row 86 ---------------------------------------------------------------
$ws.Range("A86").Value = "'6979"
$ws.Range("B86").Value = "'8/18/2025"
$ws.Range("C86").Value = "RIVADAVIA AV. 6740"
$ws.Range("E86").Value = "'809006419"
$ws.Range("H86").Value = "Reclaman fuera de plomo ver si es necesario cambio"
$ws.Range("J86").Value = "Aplomo"
$ws.Range("M86").Value = -58.460441
$ws.Range("N86").Value = -34.628243

# --- Row 87 ---------------------------------------------------------------
$ws.Range("A87").Value = "'-557"
$ws.Range("B87").Value = "'8/21/2025"
$ws.Range("C87").Value = "Av Castañares 4621"
$ws.Range("D87").Value = "'8"
$ws.Range("E87").Value = "ICD30462144"
$ws.Range("H87").Value = "Colocar columna para pedir traspaso de nodo telecom"
$ws.Range("J87").Value = "Cambio"
$ws.Range("K87").Value = "Nodo Teco"
$ws.Range("M87").Value = -58.470977
$ws.Range("N87").Value = -34.665358

# --- Row 88 ---------------------------------------------------------------
$ws.Range("A88").Value = "'7051"
$ws.Range("B88").Value = "'8/26/2025"
$ws.Range("C88").Value = "MORENO, JOSE MARIA AV. 345"
$ws.Range("D88").Value = "'6"
$ws.Range("E88").Value = "Pendiente ADM"
$ws.Range("H88").Value = "Colocar PRFV R400 para pedir traspaso de fuente"
$ws.Range("K88").Value = "Fuente Teco"
$ws.Range("L88").Value = "Terminal"
$ws.Range("M88").Value = -58.435017
$ws.Range("N88").Value = -34.622044

# --- Row 89 ---------------------------------------------------------------
$ws.Range("A89").Value = "'7060"
$ws.Range("B89").Value = "'8/26/2025"
$ws.Range("C89").Value = "PINZON 1578"
$ws.Range("D89").Value = "'4"
$ws.Range("E89").Value = "'809195671"
$ws.Range("H89").Value = "Picada"
$ws.Range("M89").Value = -58.373428
$ws.Range("N89").Value = -34.63705
$ws.Range("O89").Value = "San Telmo"
